$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Inputs sheet: secondary-inlet value 1 -> 0.1, move selection to C15
# ---------------------------------------------------------------------------
$inputs = $wb.Worksheets.Item("Inputs")
$inputs.Range("C14").Value = 0.1

# ---------------------------------------------------------------------------
# 2. Intermediate sheet: new area / bulk-flow-velocity formulas + labels
# ---------------------------------------------------------------------------
$intermediate = $wb.Worksheets.Item("Intermediate")

# ---------------------------------------------------------------------------
# 3. New "Intermediate 2" sheet - duplicate of the (still pristine)
#    "Intermediate" sheet placed right before "Outputs". Do this BEFORE the
#    "Intermediate" sheet gets its own new content below, so the copy starts
#    from the original 17-row table.
# ---------------------------------------------------------------------------
$outputs = $wb.Worksheets.Item("Outputs")
$intermediate.Copy($outputs)
$inter2 = $wb.Worksheets.Item("Intermediate (2)")
$inter2.Name = "Intermediate 2"

$intermediate.Range("A22").Value = "Area of the secondary inlet"
$intermediate.Range("A19").Value = "Bulk flow velocity"
$intermediate.Range("B22").Value = "flow rate (m2/s)"
$intermediate.Range("A23").Formula = "=(PI()/4)*((D6/100)-(D5/100))^2"
$intermediate.Range("B23").Formula = "=Inputs!C14/Intermediate!A23"

# Convert the duplicated table's C:E columns to metres (divide the original
# Intermediate-sheet cm values by 100) and relabel the unit column as "m".
for ($r = 2; $r -le 17; $r++) {
    $inter2.Range("C$r").Formula = "=Intermediate!C$r/100"
    $inter2.Range("D$r").Formula = "=Intermediate!D$r/100"
    $inter2.Range("E$r").Formula = "=Intermediate!E$r/100"
    $inter2.Range("F$r").Value = "m"
}

Write-Output "done"
